# Auto-generated COM-interop script applying the 2026-01-26 commit
# to 202601_HL_Maintain_Report.xlsx (sheet "Report").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Report")

# ---------------------------------------------------------------
# 1. Refresh the report title (reprint date 01-23 -> 01-26)
# ---------------------------------------------------------------
$ws.Range("A1").Value = "萊爾富 工作統計表  篩選月份：202601   (  製表日期:2026-01-26  )"

# ---------------------------------------------------------------
# 2. Re-wrap the last two cells of the (old) final data row 200
#    so their multi-line remarks display fully (P200, AC200).
# ---------------------------------------------------------------
$ws.Range("P200").WrapText = $true
$ws.Range("AC200").WrapText = $true

# ---------------------------------------------------------------
# 3. Append 7 new ticket rows (201-207), mirroring the banded
#    formatting of the existing table (odd rows -> row 3 style
#    template, even rows -> row 4 style template).
# ---------------------------------------------------------------

# --- Row 201 ------------------------------------------------------------
$ws.Range("A3:AK3").Copy() | Out-Null
$ws.Range("A201:AK201").PasteSpecial(-4122) | Out-Null
$ws.Range("A201").Value = 199
$ws.Range("B201").Value = "維修"
$ws.Range("C201").Value = 2026013396
$ws.Range("D201").NumberFormat = "@"
$ws.Range("D201").Value = "13569115012401"
$ws.Range("E201").Value = "急修件"
$ws.Range("F201").Value = 3569
$ws.Range("G201").Value = "中和中板店"
$ws.Range("H201").Value = "新北市中和區"
$ws.Range("I201").Value = "2026-01-24 09:47:11"
$ws.Range("J201").Value = "星期六"
$ws.Range("K201").Value = "上午"
$ws.Range("L201").Value = "HLD3"
$ws.Range("M201").Value = "HL-熱感發票機"
$ws.Range("N201").Value = "D308"
$ws.Range("O201").Value = "上蓋打不開"
$ws.Range("P201").Value = "01/24 09:51百大門市啟動緊急叫修:門市反應TM1發票機(BSC-10、BSC10II)更換紙捲後上蓋打開後無法合起來，門市告知無法自行排除...台芝到店協助"
$ws.Range("Q201").Value = "THILF03569"
$ws.Range("R201").Value = "新北一"
$ws.Range("S201").Value = "劉柏均"
$ws.Range("T201").Value = 1
$ws.Range("U201").Value = "已完工"
$ws.Range("V201").Value = "2026-01-24 09:52:37"
$ws.Range("W201").Value = "2026-01-24 11:10:00"
$ws.Range("X201").Value = "2026-01-24 11:28:00"
$ws.Range("Y201").Value = "2026-01-24 15:52:00"
$ws.Range("Z201").Value = 0.3
$ws.Range("AB201").Value = "到場處理"
$ws.Range("AC201").Value = "裁刀為定位. 清潔後即可"
$ws.Range("AK201").Value = "O"
$ws.Range("D3").Copy() | Out-Null
$ws.Range("D201").PasteSpecial(-4122) | Out-Null
$ws.Range("P201").WrapText = $true
$ws.Range("AC201").WrapText = $true

# --- Row 202 ------------------------------------------------------------
$ws.Range("A4:AK4").Copy() | Out-Null
$ws.Range("A202:AK202").PasteSpecial(-4122) | Out-Null
$ws.Range("A202").Value = 200
$ws.Range("B202").Value = "服務"
$ws.Range("C202").Value = 2026013401
$ws.Range("F202").Value = 3569
$ws.Range("G202").Value = "中和中板店"
$ws.Range("H202").Value = "新北市中和區"
$ws.Range("Q202").Value = "THILF03569"
$ws.Range("R202").Value = "新北一"
$ws.Range("S202").Value = "劉柏均"
$ws.Range("T202").Value = 1
$ws.Range("U202").Value = "已完工"
$ws.Range("V202").Value = "2026-01-24 11:33:26"
$ws.Range("W202").Value = "2026-01-24 11:05:00"
$ws.Range("X202").Value = "2026-01-24 11:25:00"
$ws.Range("Z202").Value = 0.3
$ws.Range("AB202").Value = "到場處理"
$ws.Range("AC202").Value = "PMQ1+7210002967"
$ws.Range("AD202").Value = "O"
$ws.Range("AJ202").Value = "O"
$ws.Range("AK202").Value = "O"
$ws.Range("P202").WrapText = $true
$ws.Range("AC202").WrapText = $true

# --- Row 203 ------------------------------------------------------------
$ws.Range("A3:AK3").Copy() | Out-Null
$ws.Range("A203:AK203").PasteSpecial(-4122) | Out-Null
$ws.Range("A203").Value = 201
$ws.Range("B203").Value = "維修"
$ws.Range("C203").Value = 2026013402
$ws.Range("D203").NumberFormat = "@"
$ws.Range("D203").Value = "14208115012402"
$ws.Range("E203").Value = "急修件"
$ws.Range("F203").Value = 4208
$ws.Range("G203").Value = "新莊瓊泰店"
$ws.Range("H203").Value = "新北市新莊區"
$ws.Range("I203").Value = "2026-01-24 11:53:29"
$ws.Range("J203").Value = "星期六"
$ws.Range("K203").Value = "上午"
$ws.Range("L203").Value = "HL24"
$ws.Range("M203").Value = "HL-SC主機"
$ws.Range("N203").Value = 2405
$ws.Range("O203").Value = "檔案損毀(更換硬碟)"
$ws.Range("P203").Value = "1/24 11:50 與總公司明翰確認啟動緊急叫修:SC(SHUTTLE6S)-門市反應TM1.2主檔未更新，顯示01/23，應為01/24，查看SC主檔01/24，執行SC轉TM>手動TM版更>10分鐘後TM執行版更仍異常，客服嘗試至最高權限執行手動排除主檔無法更新>查詢SC主檔為當日>SC轉TM>手動TM版本更新>10分鐘後TM執行版本更新仍異常，經HIPOS執行NewOpen->手動主檔轉入->Sc轉tm後，查看sc e檔 trans_in 未丟任何檔案，經總公司圭連告知疑似SC第二硬碟區塊毀損，需到店更換硬碟，並攜帶第一顆硬碟備著...請台芝到店協助 PS.若因更換HD.請跟店長宣達:1.請門市先回報代收會計 2.請確認SC的代收資料是否正確 (須與代收單據逐一核對) 與門市確認帳關到01/23，與通訊圭連確認都有收到"
$ws.Range("Q203").Value = "THILF04208"
$ws.Range("R203").Value = "新北一"
$ws.Range("S203").Value = "湯家瑋"
$ws.Range("T203").Value = 1
$ws.Range("U203").Value = "已完工"
$ws.Range("V203").Value = "2026-01-24 12:01:58"
$ws.Range("W203").Value = "2026-01-24 12:30:00"
$ws.Range("X203").Value = "2026-01-24 15:50:00"
$ws.Range("Y203").Value = "2026-01-24 18:01:00"
$ws.Range("Z203").Value = 3.3
$ws.Range("AB203").Value = "到場處理"
$ws.Range("AC203").Value = "更換第一二顆硬碟不備份還原
並告知李先生更換硬碟注意事項
回報0800"
$ws.Range("AK203").Value = "O"
$ws.Range("D3").Copy() | Out-Null
$ws.Range("D203").PasteSpecial(-4122) | Out-Null
$ws.Range("P203").WrapText = $true
$ws.Range("AC203").WrapText = $true

# --- Row 204 ------------------------------------------------------------
$ws.Range("A4:AK4").Copy() | Out-Null
$ws.Range("A204:AK204").PasteSpecial(-4122) | Out-Null
$ws.Range("A204").Value = 202
$ws.Range("B204").Value = "維修"
$ws.Range("C204").Value = 2026013406
$ws.Range("D204").Value = "E4155115012401"
$ws.Range("E204").Value = "一般件"
$ws.Range("F204").Value = 4155
$ws.Range("G204").Value = "三重穀保店"
$ws.Range("H204").Value = "新北市三重區"
$ws.Range("I204").Value = "2026-01-24 17:46:54"
$ws.Range("J204").Value = "星期六"
$ws.Range("K204").Value = "下午"
$ws.Range("L204").Value = "HL60"
$ws.Range("M204").Value = "HL-LIFE-ET印票機L90"
$ws.Range("N204").Value = 6004
$ws.Range("O204").Value = "無電源反應、無法開機"
$ws.Range("P204").Value = "L90門市反應昨天停電復電後票卷機無法開機無電源反應，門市已有重新拔插按電源健仍異常無反應...須請台芝到店協助(昨晚台電停電 票卷機故障 無法開機)"
$ws.Range("Q204").Value = "THILF04155"
$ws.Range("R204").Value = "新北一"
$ws.Range("S204").Value = "劉柏均"
$ws.Range("T204").Value = 1
$ws.Range("U204").Value = "已完工"
$ws.Range("V204").Value = "2026-01-24 17:53:36"
$ws.Range("W204").Value = "2026-01-26 16:45:00"
$ws.Range("X204").Value = "2026-01-26 16:56:00"
$ws.Range("Y204").Value = "2026-01-27 13:00:00"
$ws.Range("Z204").Value = 0.2
$ws.Range("AB204").Value = "到場處理"
$ws.Range("AC204").Value = "重新插拔電源. "
$ws.Range("AK204").Value = "O"
$ws.Range("P204").WrapText = $true
$ws.Range("AC204").WrapText = $true

# --- Row 205 ------------------------------------------------------------
$ws.Range("A3:AK3").Copy() | Out-Null
$ws.Range("A205:AK205").PasteSpecial(-4122) | Out-Null
$ws.Range("A205").Value = 203
$ws.Range("B205").Value = "服務"
$ws.Range("C205").Value = 2026013468
$ws.Range("F205").Value = 3627
$ws.Range("G205").Value = "三重重新橋"
$ws.Range("H205").Value = "新北市三重區"
$ws.Range("Q205").Value = "THILF03627"
$ws.Range("R205").Value = "新北一"
$ws.Range("S205").Value = "吳宗鴻"
$ws.Range("T205").Value = 1
$ws.Range("U205").Value = "已完工"
$ws.Range("V205").Value = "2026-01-26 14:14:40"
$ws.Range("W205").Value = "2026-01-26 11:00:00"
$ws.Range("X205").Value = "2026-01-26 14:13:00"
$ws.Range("Z205").Value = 3.2
$ws.Range("AB205").Value = "到場處理"
$ws.Range("AC205").Value = "裝潢回裝完成"
$ws.Range("AE205").Value = "O"
$ws.Range("AK205").Value = "O"
$ws.Range("P205").WrapText = $true
$ws.Range("AC205").WrapText = $true

# --- Row 206 ------------------------------------------------------------
$ws.Range("A4:AK4").Copy() | Out-Null
$ws.Range("A206:AK206").PasteSpecial(-4122) | Out-Null
$ws.Range("A206").Value = 204
$ws.Range("B206").Value = "服務"
$ws.Range("C206").Value = 2026013473
$ws.Range("F206").Value = 3627
$ws.Range("G206").Value = "三重重新橋"
$ws.Range("H206").Value = "新北市三重區"
$ws.Range("Q206").Value = "THILF03627"
$ws.Range("R206").Value = "新北一"
$ws.Range("S206").Value = "劉柏均"
$ws.Range("T206").Value = 1
$ws.Range("U206").Value = "已完工"
$ws.Range("V206").Value = "2026-01-26 14:23:46"
$ws.Range("W206").Value = "2026-01-26 11:00:00"
$ws.Range("X206").Value = "2026-01-26 14:00:00"
$ws.Range("Z206").Value = 3
$ws.Range("AC206").Value = "回裝已完工已請0800確認版本"
$ws.Range("AE206").Value = "O"
$ws.Range("AK206").Value = "O"
$ws.Range("P206").WrapText = $true
$ws.Range("AC206").WrapText = $true

# --- Row 207 ------------------------------------------------------------
$ws.Range("A3:AK3").Copy() | Out-Null
$ws.Range("A207:AK207").PasteSpecial(-4122) | Out-Null
$ws.Range("A207").Value = 205
$ws.Range("B207").Value = "服務"
$ws.Range("C207").Value = 2026013499
$ws.Range("F207").Value = 4155
$ws.Range("G207").Value = "三重穀保店"
$ws.Range("H207").Value = "新北市三重區"
$ws.Range("Q207").Value = "THILF04155"
$ws.Range("R207").Value = "新北一"
$ws.Range("S207").Value = "劉柏均"
$ws.Range("T207").Value = 1
$ws.Range("U207").Value = "已完工"
$ws.Range("V207").Value = "2026-01-26 16:59:22"
$ws.Range("W207").Value = "2026-01-26 16:38:00"
$ws.Range("X207").Value = "2026-01-26 16:55:00"
$ws.Range("Z207").Value = 0.3
$ws.Range("AB207").Value = "到場處理"
$ws.Range("AC207").Value = "PMQ1+7210002791"
$ws.Range("AD207").Value = "O"
$ws.Range("AJ207").Value = "O"
$ws.Range("AK207").Value = "O"

# ---------------------------------------------------------------
# 4. Extend the print area to cover the newly added rows
# ---------------------------------------------------------------
$wb.Names.Item("Report!Print_Area").RefersTo = "='Report'!`$A`$1:`$AK`$207"

# ---------------------------------------------------------------
# 5. Restore the cursor/selection to the new last row, as saved
#    by the author after entering the data
# ---------------------------------------------------------------
$ws.Range("A207").Select() | Out-Null

